$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New StatQuery text (replaces the old, buggy Cypher stat query used for the
# faceted-filter counts) - commit: "Fixed StudyComb for Faceted Filters ICDC"
$newStatQuery = "MATCH (demo:demographic)`nWHERE demo.breed IN [`"American Staffordshire Terrier`"]`nMATCH (demo:demographic)-->(c:case)-->(s:study)-->(p:program)`nOPTIONAL MATCH (c)<-[*]-(samp:sample)`nOPTIONAL MATCH (c)<-[*]-(f:file)`nRETURN `n`tcount(DISTINCT(f)) as number_of_files, `n`tcount(DISTINCT(samp)) as number_of_sample, `n`tcount(DISTINCT(c)) as number_of_cases, `n`tcount(DISTINCT(s)) as number_of_study"

# Column C (StatQuery) holds the same query for the three tab rows (Cases,
# Samples, Files) - update all three to the corrected query.
$ws.Range("C2").Value = $newStatQuery
$ws.Range("C3").Value = $newStatQuery
$ws.Range("C4").Value = $newStatQuery

# Narrow column C a bit (was 123.33, now ~94.11 characters wide)
$ws.Columns.Item(3).ColumnWidth = 93.33

# Scroll/zoom the view back to normal (was zoomed to 145%, showing row 4 at
# top); reset to 100% zoom, showing row 3 at top.
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 1
$win.Zoom = 100
